$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.618.91'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.594.89'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.02'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.244'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.36'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0837'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D13').Value = '1.588.06'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.55'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '26.592.91'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '207.37'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.88'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.25'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.30'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.87'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.34'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.114'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0504'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.16'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.654'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = '1.283.88'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.786'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.59'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.919'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +9.26%  '
$ws.Range('D46').Value = '1.730.81'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '89.66'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('E51').Value = '  -1.46%  '
